# feat: add 2022-Q1 data
#
# The existing "总计" (grand-total) sheet is logically split in two:
#   1. It is renamed to "2022-Q1" and repopulated with the new quarter's
#      per-fund holding detail (mirrors the layout of the other quarterly
#      sheets: 2020-Q4 .. 2021-Q4).
#   2. A duplicate of it is appended right after, renamed back to "总计",
#      and reproduces the grand-total-by-quarter table with a new leading
#      row for "2022-Q1" (the older rows shift down by one).
#
# Both new sheets are produced via Worksheet.Copy() off the original
# "总计" sheet (rather than Worksheets.Add()) so they inherit its exact
# page setup / header+index-column cell styling instead of picking up
# the blank-workbook defaults.

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

$original = $wb.Worksheets.Item("总计")
$original.Copy($null, $original)

$q1 = $wb.Worksheets.Item("总计")
$total = $wb.Worksheets.Item("总计 (2)")

$q1.Name = "2022-Q1"
$total.Name = "总计"

# ---------------------------------------------------------------------
# Step 1: "2022-Q1" — replace the grand-total table with the fund-level
# holding detail table (A1:H8).
# ---------------------------------------------------------------------

# extend the styled header row (currently B1:D1) and index column
# (currently A2:A6) out to the new table's extent before touching values,
# so the added cells inherit the same styling (bold / border / center-top)
# as the existing ones.
$q1.Range("B1").Copy() | Out-Null
$q1.Range("E1:H1").PasteSpecial($xlPasteFormats) | Out-Null
$q1.Range("A2").Copy() | Out-Null
$q1.Range("A7:A8").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q1.Cells.Item(1, $col).Value = $headers[$col - 2]
}

$fundRows = @(
    @("008415", "国泰大制造两年持有期混合", "23.19", "92.05", "2.85", "0.6609", 8),
    @("020010", "国泰金牛创新混合",         "16.99", "84.21", "2.74", "0.4655", 10),
    @("011466", "兴业医疗保健混合A",        "7.10",  "84.68", "4.37", "0.3103", 7),
    @("020026", "国泰成长优选混合",         "8.72",  "92.38", "3.33", "0.2904", 6),
    @("011467", "兴业医疗保健混合C",        "2.22",  "84.68", "4.37", "0.0970", 7),
    @("008618", "永赢医药健康股票A",        "0.77",  "90.10", "6.65", "0.0512", 4),
    @("008619", "永赢医药健康股票C",        "0.39",  "90.10", "6.65", "0.0259", 4)
)

# the fund-code / scale / position / ratio / value columns must stay as
# literal text (not get auto-converted to numbers) so leading zeros in
# codes like "008415" and trailing zeros like "90.10" are preserved, same
# as every other quarterly sheet in this workbook.
$q1.Range("B2:G8").NumberFormat = "@"

$r = 2
foreach ($row in $fundRows) {
    $q1.Cells.Item($r, 1).Value = $r - 2
    $q1.Cells.Item($r, 2).Value = $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = $row[2]
    $q1.Cells.Item($r, 5).Value = $row[3]
    $q1.Cells.Item($r, 6).Value = $row[4]
    $q1.Cells.Item($r, 7).Value = $row[5]
    $q1.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# ---------------------------------------------------------------------
# Step 2: "总计" — keep the existing 3-column layout (日期 / 持有数量(只) /
# 持有市值(亿元)) but grow from 5 to 6 data rows: a new "2022-Q1" row up
# top, followed by the previously-existing rows (re-indexed in column A).
# ---------------------------------------------------------------------
$total.Range("A2").Copy() | Out-Null
$total.Range("A7").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

$totalRows = @(
    @("2022-Q1", 7, 1.9),
    @("2021-Q4", 5, 1),
    @("2021-Q3", 1, 0.79),
    @("2021-Q2", 6, 1.89),
    @("2021-Q1", 7, 4.15),
    @("2020-Q4", 4, 0.15)
)

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $r - 2
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
    $r++
}

Write-Output "2022-Q1 sheet added; 总计 sheet rebuilt."
